$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("D2") "28.398.77"
Set-TextValue $ws.Range("E2") "  -0.06%  "
Set-TextValue $ws.Range("D3") "1.821.61"
Set-TextValue $ws.Range("E3") "  -0.15%  "
Set-TextValue $ws.Range("D4") "1.004"
Set-TextValue $ws.Range("E4") "  +0.09%  "
Set-TextValue $ws.Range("D5") "315.55"
Set-TextValue $ws.Range("E5") "  +0.23%  "
Set-TextValue $ws.Range("D6") "1.003"
Set-TextValue $ws.Range("E6") "  +0.10%  "
Set-TextValue $ws.Range("D7") "0.5229"
Set-TextValue $ws.Range("E7") "  +2.40%  "
Set-TextValue $ws.Range("E8") "  -1.63%  "
Set-TextValue $ws.Range("E9") "  +5.36%  "
Set-TextValue $ws.Range("B10") "Polygon"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D10") "1.115"
Set-TextValue $ws.Range("E10") "  +0.88%  "
Set-TextValue $ws.Range("B11") "OKB"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D11") "41.87"
Set-TextValue $ws.Range("E11") "  +0.70%  "
Set-TextValue $ws.Range("D12") "6.395"
Set-TextValue $ws.Range("E12") "  +2.14%  "
Set-TextValue $ws.Range("E13") "  +0.09%  "
Set-TextValue $ws.Range("D14") "20.88"
Set-TextValue $ws.Range("E14") "  -0.57%  "
Set-TextValue $ws.Range("D15") "7.431"
Set-TextValue $ws.Range("E15") "  -0.94%  "
Set-TextValue $ws.Range("D16") "1.816.70"
Set-TextValue $ws.Range("E16") "  -0.34%  "
Set-TextValue $ws.Range("D17") "94.36"
Set-TextValue $ws.Range("E17") "  +1.17%  "
Set-TextValue $ws.Range("D18") "0.00001105"
Set-TextValue $ws.Range("E18") "  +1.01%  "
Set-TextValue $ws.Range("E19") "  -0.46%  "
Set-TextValue $ws.Range("D20") "17.61"
Set-TextValue $ws.Range("E20") "  +0.01%  "
Set-TextValue $ws.Range("D21") "1.003"
Set-TextValue $ws.Range("E21") "  +0.15%  "
Set-TextValue $ws.Range("D22") "6.021"
Set-TextValue $ws.Range("E22") "  -1.89%  "
Set-TextValue $ws.Range("D23") "28.444.90"
Set-TextValue $ws.Range("E23") "  +0.00%  "
Set-TextValue $ws.Range("D24") "11.36"
Set-TextValue $ws.Range("E24") "  +1.85%  "
Set-TextValue $ws.Range("D25") "2.245"
Set-TextValue $ws.Range("E25") "  -0.80%  "
Set-TextValue $ws.Range("D26") "159.11"
Set-TextValue $ws.Range("E26") "  +1.90%  "
Set-TextValue $ws.Range("D27") "20.86"
Set-TextValue $ws.Range("E27") "  +0.82%  "
Set-TextValue $ws.Range("D28") "2.032.28"
Set-TextValue $ws.Range("E28") "  -0.05%  "
Set-TextValue $ws.Range("D29") "2.412"
Set-TextValue $ws.Range("E29") "  +1.28%  "
Set-TextValue $ws.Range("D30") "124.47"
Set-TextValue $ws.Range("E30") "  +0.40%  "
Set-TextValue $ws.Range("D31") "0.1110"
Set-TextValue $ws.Range("E31") "  +2.14%  "
Set-TextValue $ws.Range("D32") "1.077"
Set-TextValue $ws.Range("E32") "  -2.62%  "
Set-TextValue $ws.Range("D33") "5.676"
Set-TextValue $ws.Range("E33") "  +0.89%  "
Set-TextValue $ws.Range("D34") "3.683"
Set-TextValue $ws.Range("E34") "  +0.73%  "
Set-TextValue $ws.Range("D35") "0.07307"
Set-TextValue $ws.Range("E35") "  +3.88%  "
Set-TextValue $ws.Range("D36") "12.23"
Set-TextValue $ws.Range("E36") "  +9.11%  "
Set-TextValue $ws.Range("D37") "0.2205"
Set-TextValue $ws.Range("E37") "  +0.08%  "
Set-TextValue $ws.Range("D38") "0.02348"
Set-TextValue $ws.Range("E38") "  +1.46%  "
Set-TextValue $ws.Range("B39") "FraxShare"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D39") "8.772"
Set-TextValue $ws.Range("E39") "  -0.48%  "
Set-TextValue $ws.Range("B40") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D40") "5.118"
Set-TextValue $ws.Range("E40") "  -0.55%  "
Set-TextValue $ws.Range("D41") "0.6312"
Set-TextValue $ws.Range("E41") "  +1.32%  "
Set-TextValue $ws.Range("D42") "1.182"
Set-TextValue $ws.Range("E42") "  +1.05%  "
Set-TextValue $ws.Range("E43") "  -0.33%  "
Set-TextValue $ws.Range("B44") "Decentraland"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D44") "0.6143"
Set-TextValue $ws.Range("E44") "  +4.56%  "
Set-TextValue $ws.Range("B45") "EnergySwap"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "13.42"
Set-TextValue $ws.Range("E45") "  +0.29%  "
Set-TextValue $ws.Range("D46") "3.799"
Set-TextValue $ws.Range("E46") "  +2.50%  "
Set-TextValue $ws.Range("D47") "127.13"
Set-TextValue $ws.Range("E47") "  +1.84%  "
Set-TextValue $ws.Range("E48") "  +2.22%  "
Set-TextValue $ws.Range("D49") "1.970"
Set-TextValue $ws.Range("E49") "  -0.13%  "
Set-TextValue $ws.Range("D50") "0.06902"
Set-TextValue $ws.Range("E50") "  -0.09%  "
Set-TextValue $ws.Range("D51") "73.99"
Set-TextValue $ws.Range("E51") "  +0.18%  "
